$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = "0"
$ws.Range("G15").Value = "0"
$ws.Range("H15").Value = "***.*"
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 0

# --- Row 16 (Fel. Assault) ---
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 11
$ws.Range("H16").Value = -26.666666666666
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 47
$ws.Range("K16").Value = -6.382978723404
$ws.Range("L16").Value = 57.142857142857
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = -76.719576719576

# --- Row 17 (Burglary) ---
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 36
$ws.Range("K17").Value = 11.111111111111
$ws.Range("L17").Value = 37.931034482758
$ws.Range("M17").Value = -6.976744186046
$ws.Range("N17").Value = -45.945945945945

# --- Row 18 (Gr. Larceny) ---
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -23.529411764705
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 44
$ws.Range("K18").Value = -27.272727272727
$ws.Range("L18").Value = -15.789473684210
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = -78.231292517006

# --- Row 19 (G.L.A.) ---
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 18.75
$ws.Range("F19").Value = 61
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 45.238095238095
$ws.Range("I19").Value = 222
$ws.Range("J19").Value = 196
$ws.Range("K19").Value = 13.265306122449
$ws.Range("L19").Value = 76.190476190476
$ws.Range("M19").Value = 20
$ws.Range("N19").Value = -13.618677042801

# --- Row 20 (TOTAL-ish row, no C/D/E edits) ---
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 166.666666666667
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 18
$ws.Range("K20").Value = -5.555555555555
$ws.Range("L20").Value = 183.333333333333
$ws.Range("M20").Value = 6.25
$ws.Range("N20").Value = -90.502793296089

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 3.448275862068
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = 19.540229885057
$ws.Range("I21").Value = 357
$ws.Range("J21").Value = 346
$ws.Range("K21").Value = 3.179190751445
$ws.Range("L21").Value = 55.217391304347
$ws.Range("M21").Value = 18.604651162790
$ws.Range("N21").Value = -58.245614035087

# --- Row 22 (Transit) ---
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 13
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = 44.444444444444
$ws.Range("L22").Value = 333.333333333333
$ws.Range("M22").Value = 333.333333333333

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 300
$ws.Range("I23").Value = 15
$ws.Range("K23").Value = -28.571428571428
$ws.Range("L23").Value = -37.5
$ws.Range("M23").Value = -6.25

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 8
$ws.Range("E24").Value = -20
$ws.Range("F24").Value = 55
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = 30.952380952381
$ws.Range("I24").Value = 202
$ws.Range("J24").Value = 190
$ws.Range("K24").Value = 6.315789473684
$ws.Range("L24").Value = 42.253521126760
$ws.Range("M24").Value = -16.528925619834

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 4.545454545454
$ws.Range("I25").Value = 106
$ws.Range("J25").Value = 78
$ws.Range("K25").Value = 35.897435897435
$ws.Range("L25").Value = 51.428571428571
$ws.Range("M25").Value = 4.950495049504

# --- Row 26 (UCR Rape*) ---
$ws.Range("F26").Value = "0"
$ws.Range("G26").Value = "0"
$ws.Range("H26").Value = "***.*"
$ws.Range("L26").Value = 0

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 23
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 43.75
$ws.Range("L27").Value = 228.571428571429

# --- Row 30 (Hate Crimes) ---
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 9
$ws.Range("K30").Value = 50
$ws.Range("L30").Value = 800
